$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update values in column A (values scraped from updated KNN imputation run)
$ws.Cells.Item(8, 1).Value = -22.305
$ws.Cells.Item(10, 1).Value = -21.846
$ws.Cells.Item(12, 1).Value = -21.546
$ws.Cells.Item(18, 1).Value = -22.095
$ws.Cells.Item(37, 1).Value = -19.997
$ws.Cells.Item(55, 1).Value = -22.184
$ws.Cells.Item(68, 1).Value = -21.567
$ws.Cells.Item(77, 1).Value = -20.666
$ws.Cells.Item(78, 1).Value = -19.951
$ws.Cells.Item(81, 1).Value = -21.776
$ws.Cells.Item(82, 1).Value = -22.067
